$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label cell for the new averages row, italic
$ws.Range("A1").Value = "Avrg"
$ws.Range("A1").Font.Italic = $true

# Also italicize the existing header label in A2 (" Temp C")
$ws.Range("A2").Font.Italic = $true

# Average formulas for each data column (rows 3 through 69, since row 2 is the header row)
$cols = @("B","C","D","E","F","G","H","I","J","K","L")
foreach ($col in $cols) {
    $ws.Range("$col" + "1").Formula = "=AVERAGE(" + $col + "3:" + $col + "69)"
}

# Reset the cursor to the default top-left cell (clears the stale M2 selection
# that was left over in the source file)
$ws.Range("A1").Select()
